$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds the new "Rogue" class data, mirroring the rows already
# present for "Cleric" in column A. Copy column A's formatting across first
# so the new cells pick up the same style, then overwrite with Rogue's
# values. Row 2 is a generic (class-independent) proficiency-bonus-by-level
# table, so it stays identical in both columns.
$ws.Range("A1:A9").Copy($ws.Range("B1:B9"))

$ws.Range("B1").Value = "Rogue"
$ws.Range("B3").Value = 8
$ws.Range("B4").Value = "LA"
$ws.Range("B5").Value = "S=HC=LS=R=SS"
$ws.Range("B6").Value = "Thieves' Tools"
$ws.Range("B7").Value = "Dexterity=Intelligence"
$ws.Range("B8").Value = "false=None"
$ws.Range("B9").Value = "1/Thieves' Cant=1/Sneak Attack=2/Cunning Action=5/Uncanny Dodge=7/Evasion=11/Reliable Talent=14/Blindsense=15/Slippery Mind=18/Elusive=20/Stroke of Luck"
